# Updates the 20x5 math-problem table (100 cells) to the new answer set.
# Cell text is replaced positionally (row-major order) via Cell.Range.Text,
# since several old values are duplicated (e.g. "6+58=64" appears twice but
# maps to two different new values), ruling out a global text Find/Replace.

$d = $word.ActiveDocument

$newValues = @(
    "17+15=32",
    "9+63=72",
    "92-44=48",
    "46-37=9",
    "69+23=92",
    "58+37=95",
    "9+34=43",
    "72-3=69",
    "12+9=21",
    "28+23=51",
    "6+48=54",
    "5+16=21",
    "46+27=73",
    "25+66=91",
    "23+38=61",
    "39+42=81",
    "49+38=87",
    "16+58=74",
    "67-8=59",
    "81-2=79",
    "40-22=18",
    "49+26=75",
    "29+25=54",
    "79+2=81",
    "72-55=17",
    "16+39=55",
    "9+67=76",
    "9+18=27",
    "97-78=19",
    "93-44=49",
    "86-37=49",
    "52+19=71",
    "72-68=4",
    "25+67=92",
    "28+27=55",
    "15+27=42",
    "8+78=86",
    "63-54=9",
    "78+15=93",
    "95-36=59",
    "24+17=41",
    "62-6=56",
    "28+26=54",
    "20-17=3",
    "55-47=8",
    "43-34=9",
    "70-56=14",
    "27+28=55",
    "92-3=89",
    "84-65=19",
    "58+16=74",
    "53-18=35",
    "70-53=17",
    "54-19=35",
    "37+16=53",
    "61-24=37",
    "18+74=92",
    "56+37=93",
    "69+5=74",
    "71-38=33",
    "16+37=53",
    "71-43=28",
    "40-34=6",
    "36+39=75",
    "54-5=49",
    "47+19=66",
    "55-6=49",
    "62-34=28",
    "96-59=37",
    "29+26=55",
    "49+6=55",
    "6+18=24",
    "52-6=46",
    "25+16=41",
    "13+59=72",
    "19+69=88",
    "38+39=77",
    "51-42=9",
    "80-38=42",
    "35-28=7",
    "28+48=76",
    "91-29=62",
    "76-38=38",
    "53-25=28",
    "9+37=46",
    "64-38=26",
    "19+74=93",
    "70-59=11",
    "39+22=61",
    "95-7=88",
    "46+28=74",
    "55+7=62",
    "8+47=55",
    "71-37=34",
    "20-15=5",
    "36-29=7",
    "19+48=67",
    "9+87=96",
    "17+49=66",
    "48+45=93"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output ("Updated " + $idx + " cells")
